# Auto-generated edit script: updates cached market-price / profit
# columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to
# reflect a refreshed data pull from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 506.85
$ws.Range("I19").Value = 522.5
$ws.Range("K19").Value = 522.5
$ws.Range("M19").Value = -347.5
$ws.Range("H100").Value = 4301.1924
$ws.Range("I100").Value = 2137.0715
$ws.Range("J100").Value = 6826
$ws.Range("K100").Value = 2137.0715
$ws.Range("L100").Value = 6826
$ws.Range("M100").Value = -1596.0715
$ws.Range("N100").Value = -7908
$ws.Range("H113").Value = 4372.1577
$ws.Range("J113").Value = 4532.7827
$ws.Range("L113").Value = 4532.7827
$ws.Range("N113").Value = -11040.7827
$ws.Range("H116").Value = 4169.5
$ws.Range("I116").Value = 3551.6365
$ws.Range("K116").Value = 3551.6365
$ws.Range("M116").Value = -109.6365000000001
$ws.Range("H137").Value = 3263.1428
$ws.Range("I137").Value = 2828.6086
$ws.Range("K137").Value = 8485.825800000001
$ws.Range("M137").Value = -5935.825800000001
$ws.Range("H138").Value = 2931.8406
$ws.Range("I138").Value = 1373.12
$ws.Range("J138").Value = 3817.4773
$ws.Range("K138").Value = 4119.36
$ws.Range("L138").Value = 11452.4319
$ws.Range("M138").Value = 1020.64
$ws.Range("N138").Value = -21732.4319

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4817.706
$ws.Range("I45").Value = 3754.625
$ws.Range("K45").Value = 3754.625
$ws.Range("M45").Value = -3377.625
$ws.Range("H80").Value = 25364
$ws.Range("I80").Value = 25364
$ws.Range("K80").Value = 25364
$ws.Range("M80").Value = -24366
$ws.Range("H83").Value = 25364
$ws.Range("I83").Value = 25364
$ws.Range("K83").Value = 76092
$ws.Range("M83").Value = -71100
$ws.Range("H102").Value = 1141.9412
$ws.Range("I102").Value = 1141.9412
$ws.Range("K102").Value = 1141.9412
$ws.Range("M102").Value = 480.0588
$ws.Range("H122").Value = 4021.0264
$ws.Range("I122").Value = 2026.6207
$ws.Range("K122").Value = 6079.8621
$ws.Range("M122").Value = -3629.8621
$ws.Range("H132").Value = 3249.3962
$ws.Range("I132").Value = 2964.2727
$ws.Range("J132").Value = 4643.3335
$ws.Range("K132").Value = 8892.8181
$ws.Range("L132").Value = 13930.0005
$ws.Range("M132").Value = -6362.8181
$ws.Range("N132").Value = -18990.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 14250
$ws.Range("J95").Value = 14250
$ws.Range("L95").Value = 14250
$ws.Range("N95").Value = -19742
$ws.Range("H99").Value = 35369.918
$ws.Range("I99").Value = 38130.816
$ws.Range("K99").Value = 38130.816
$ws.Range("M99").Value = -36632.816

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 619.75
$ws.Range("I16").Value = 381.1111
$ws.Range("K16").Value = 381.1111
$ws.Range("M16").Value = -94.11110000000002
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H22").Value = 333498.34
$ws.Range("I22").Value = 215.4
$ws.Range("J22").Value = 500139.8
$ws.Range("K22").Value = 215.4
$ws.Range("L22").Value = 500139.8
$ws.Range("M22").Value = 134.6
$ws.Range("N22").Value = -500839.8
$ws.Range("H31").Value = 4326.636
$ws.Range("I31").Value = 1723.8334
$ws.Range("K31").Value = 1723.8334
$ws.Range("M31").Value = -1428.8334
$ws.Range("H34").Value = 4326.636
$ws.Range("I34").Value = 1723.8334
$ws.Range("K34").Value = 1723.8334
$ws.Range("M34").Value = -1521.8334
$ws.Range("H107").Value = 3733.6875
$ws.Range("I107").Value = 510.5909
$ws.Range("J107").Value = 10824.5
$ws.Range("K107").Value = 510.5909
$ws.Range("L107").Value = 10824.5
$ws.Range("M107").Value = 1409.4091
$ws.Range("N107").Value = -14664.5
$ws.Range("H113").Value = 619.75
$ws.Range("I113").Value = 381.1111
$ws.Range("K113").Value = 381.1111
$ws.Range("M113").Value = 1788.8889
$ws.Range("H132").Value = 3003.0833
$ws.Range("I132").Value = 2440.0938
$ws.Range("K132").Value = 7320.2814
$ws.Range("M132").Value = -4790.2814

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 1971.4286
$ws.Range("I115").Value = 1950
$ws.Range("K115").Value = 5850
$ws.Range("M115").Value = -4675
$ws.Range("H117").Value = 916.4167
$ws.Range("I117").Value = 507.83334
$ws.Range("K117").Value = 1523.50002
$ws.Range("M117").Value = 1918.49998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2409.4666
$ws.Range("I102").Value = 1843.9259
$ws.Range("K102").Value = 1843.9259
$ws.Range("M102").Value = -221.9259
$ws.Range("H126").Value = 3844.625
$ws.Range("I126").Value = 2344.3845
$ws.Range("J126").Value = 5617.636
$ws.Range("K126").Value = 7033.1535
$ws.Range("L126").Value = 16852.908
$ws.Range("M126").Value = -4563.1535
$ws.Range("N126").Value = -21792.908
$ws.Range("H132").Value = 5269.8823
$ws.Range("I132").Value = 5125
$ws.Range("J132").Value = 5535.5
$ws.Range("K132").Value = 15375
$ws.Range("L132").Value = 16606.5
$ws.Range("M132").Value = -12845
$ws.Range("N132").Value = -21666.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2726.2812
$ws.Range("I7").Value = 1465.2727
$ws.Range("K7").Value = 1465.2727
$ws.Range("M7").Value = -1353.2727
$ws.Range("H40").Value = 7619.879
$ws.Range("I40").Value = 7565.684
$ws.Range("K40").Value = 7565.684
$ws.Range("M40").Value = -7429.684
$ws.Range("H97").Value = 30336
$ws.Range("J97").Value = 30336
$ws.Range("L97").Value = 30336
$ws.Range("N97").Value = -32318
$ws.Range("H122").Value = 4977.3335
$ws.Range("I122").Value = 4021.3684
$ws.Range("K122").Value = 12064.1052
$ws.Range("M122").Value = -9614.1052
$ws.Range("H126").Value = 2726.2812
$ws.Range("I126").Value = 1465.2727
$ws.Range("K126").Value = 4395.8181
$ws.Range("M126").Value = -1925.8181
$ws.Range("H132").Value = 3326.4905
$ws.Range("I132").Value = 2875.6758
$ws.Range("K132").Value = 8627.027399999999
$ws.Range("M132").Value = -6097.027399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 461.80768
$ws.Range("I113").Value = 407.57895
$ws.Range("J113").Value = 609
$ws.Range("K113").Value = 1222.73685
$ws.Range("L113").Value = 1827
$ws.Range("M113").Value = 947.26315
$ws.Range("N113").Value = -6167
$ws.Range("H126").Value = 1809.5416
$ws.Range("I126").Value = 1736.25
$ws.Range("K126").Value = 5208.75
$ws.Range("M126").Value = -2738.75
$ws.Range("H132").Value = 2160.5454
$ws.Range("I132").Value = 2008.6842
$ws.Range("K132").Value = 6026.0526
$ws.Range("M132").Value = -3496.0526
